$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Layout:
#        A(col B)      C(col C)
# A Lag  -0.27***      -0.012*
# C Lag  -2.821***     -0.467***
#
# Update to monthly diffs/rates for 12 lags:
# Write order mirrors original shared-string table order (B2, B3, C2, C3)
$ws.Range("B2").Value = "-0.372***"
$ws.Range("B3").Value = "-3.464***"
$ws.Range("C2").Value = "0.01*"
$ws.Range("C3").Value = "-0.808***"
